$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "as of" date.
$ws.Name = "Through 2022-07-01"

# Copy the bold / centered / bordered label style (used by every row label in
# column A, cellXfs index 1) onto row 9 before anything else is edited, so the
# new "Total" row picks up the existing shared style instead of Excel minting
# a brand-new (duplicate) style entry.
$ws.Range("A1").Copy()
$ws.Range("A9").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# "June (through 06-30)" -> "June" (the month is no longer the partial one).
$ws.Range("A7").Value = "June"

# The old row 8 was the running "Total" (Jan..Jun). It now becomes the new
# "July (through 07-01)" partial-month row with the newly observed counts.
# The 2019 and 2020 columns (F and G) have no July data yet, so they are
# left blank rather than zero.
$ws.Range("A8").Value = "July (through 07-01)"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 2
$ws.Range("F8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 1

# New grand "Total" row = old Total (Jan..Jun) + the new July row above.
$ws.Range("A9").Value = "Total"
$ws.Range("B9").Value = 128
$ws.Range("C9").Value = 250
$ws.Range("D9").Value = 391
$ws.Range("E9").Value = 355
$ws.Range("F9").Value = 251
$ws.Range("G9").Value = 472
$ws.Range("H9").Value = 763
$ws.Range("I9").Value = 807
